$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 118, shifting all existing rows (118..243) down by one.
$ws.Rows.Item(118).Insert()

# Populate the newly inserted row 118 with the new weekly record.
$ws.Range("A118").Value = 8
$ws.Range("B118").Value = "Terminal La Palmera de La Serena"
$ws.Range("C118").Value = "Coquimbo"
$ws.Range("D118").Value = 44539
$ws.Range("E118").Value = 4
$ws.Range("F118").Value = 100114013
$ws.Range("G118").Value = "Zanahoria"
$ws.Range("H118").Value = "Sin especificar"
$ws.Range("I118").Value = "Primera"
$ws.Range("J118").Value = 600
$ws.Range("K118").Value = 6000
$ws.Range("L118").Value = 7000
$ws.Range("M118").Value = 6500
$ws.Range("N118").Value = "`$/saco 20 kilos"
$ws.Range("O118").Value = "Provincia del Elquí"
$ws.Range("P118").Value = 325
$ws.Range("Q118").Value = 20
$ws.Range("R118").Value = "Hortaliza"
